# fix: make most classes inherit from NamedEntity to be addressable.
#
# Appends the common NamedEntity columns (id, name, description) to every
# class sheet that doesn't already have them. ReferenceGenome and
# ReferenceSequence additionally had a stray leading "name" column that is
# dropped (its role is now played by the appended NamedEntity "name"
# column), with the remaining columns shifted left.

$wb = $excel.ActiveWorkbook

function Add-NamedEntityColumns($sheetName) {
    $ws = $wb.Worksheets.Item($sheetName)
    $lastCol = $ws.UsedRange.Columns.Count
    $idCol = $lastCol + 1
    $nameCol = $lastCol + 2
    $descCol = $lastCol + 3
    $ws.Cells.Item(1, $idCol).Value = "id"
    $ws.Cells.Item(1, $nameCol).Value = "name"
    $ws.Cells.Item(1, $descCol).Value = "description"
}

# Assay: has_sample, has_data, omics_type -> + id, name, description
Add-NamedEntityColumns "Assay"

# Sample: taxon_id, collector -> + id, name, description
Add-NamedEntityColumns "Sample"

# DataEntity: location, data_format, has_sample, has_reference -> + id, name, description
Add-NamedEntityColumns "DataEntity"

# ReferenceGenome: drop the stray leading "name" column, shift the rest
# left, then append id, name, description.
$ws = $wb.Worksheets.Item("ReferenceGenome")
$ws.Cells.Item(1, 1).Value = "has_sequence"
$ws.Cells.Item(1, 2).Value = "taxon_id"
$ws.Cells.Item(1, 3).Value = "source_uri"
$ws.Cells.Item(1, 4).Value = "id"
$ws.Cells.Item(1, 5).Value = "name"
$ws.Cells.Item(1, 6).Value = "description"

# ReferenceSequence: drop the stray leading "name" column, shift the rest
# left, then append id, name, description.
$ws = $wb.Worksheets.Item("ReferenceSequence")
$ws.Cells.Item(1, 1).Value = "location"
$ws.Cells.Item(1, 2).Value = "sequence_md5"
$ws.Cells.Item(1, 3).Value = "source_uri"
$ws.Cells.Item(1, 4).Value = "id"
$ws.Cells.Item(1, 5).Value = "name"
$ws.Cells.Item(1, 6).Value = "description"

# AlignmentSet: location, data_format, has_sample, has_reference -> + id, name, description
Add-NamedEntityColumns "AlignmentSet"

# VariantSet: location, data_format, has_sample, has_reference -> + id, name, description
Add-NamedEntityColumns "VariantSet"

# Array: location, data_format, has_sample, has_reference -> + id, name, description
Add-NamedEntityColumns "Array"
